$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 18: 경기대학교 / KGU학생부종합전형 (AI컴퓨터공학부 컴퓨터공학전공 / 20명 / 1차: 서류평가100 2차: 성적70 면접30)
$ws.Range("B18").Value = "경기대학교"
$ws.Range("E18").Value = "KGU학생부종합전형"
$ws.Range("H18").Value = "AI컴퓨터공학부 컴퓨터공학전공"
$ws.Range("K18").Value = 20
$ws.Range("L18").Value = "1차: 서류평가100`n2차: 성적70`n면접30"
$ws.Range("N18").Value = "1차 성적이 50점 미만이거나 2차 성적이 `n15점 미만인 경우 선발하지 않음"
$ws.Range("R18").Value = 3.14

# Row 20: 을지대학교(성남) / 특성화고특별전형 (의료공학전공 / 1명 / 해당없음)
$ws.Range("B20").Value = "을지대학교(성남)"
$ws.Range("E20").Value = "특성화고특별전형"
$ws.Range("H20").Value = "의료공학전공"
$ws.Range("K20").Value = 1
$ws.Range("L20").Value = "해당없음"
$ws.Range("N20").Value = "서류평가 100%(학생부 교과/비교과)"
$ws.Range("Q20").Value = "83.72`n(3.06)"
$ws.Range("R20").Value = 2.61

# Row 22: 강남대학교 (entry started, other columns left blank)
$ws.Range("B22").Value = "강남대학교"

# Remove the reference-table picture that used to float over N18:Q19 (superseded by the
# real data typed into that range).
if ($ws.Shapes.Count -gt 0) {
    for ($i = $ws.Shapes.Count; $i -ge 1; $i--) {
        $ws.Shapes.Item($i).Delete()
    }
}

# Widen column M a little so the newly typed multi-line text fits.
$ws.Columns("M").ColumnWidth = 13.71

# Leave the selection where the user was last working.
$ws.Range("E22:G23").Select()
